# Updated capital structure database
# Refreshes the per-company financial-statistics rows for the New Zealand
# "Financial Svcs. (Non-bank & Insurance)" sheet. Company names in rows 3-6 were
# re-shuffled (NZX -> General Capital -> Geneva Finance -> Blackwell Global -> NZX)
# and every metric column (D:AQ) was recomputed; some columns are sparse per row,
# so cells that no longer carry a value are cleared instead of left with stale data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: New Zealand industry aggregate (row 2) ---
$ws.Range("D2").Value = 0.1344
$ws.Range("E2").Value = 0.01325
$ws.Range("F2").Value = 0.107
$ws.Range("G2").Value = 0.2607537390146086
$ws.Range("H2").Value = 0.2607537390146086
$ws.Range("I2").Value = 0.2745081005950571
$ws.Range("J2").Value = 0.2334722844934188
$ws.Range("K2").Value = 14.335
$ws.Range("L2").Value = 0.2075460771112945
$ws.Range("M2").Value = 10.92
$ws.Range("N2").Value = 0.02532643736809147
$ws.Range("O2").Value = 0.7617718869898851
$ws.Range("P2").Value = 10.92
$ws.Range("Q2").Value = 0.02532643736809147
$ws.Range("R2").Value = 0.7617718869898851
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 36.83
$ws.Range("V2").Value = 0.08541874434677736
$ws.Range("W2").Value = 0.1072853111653447
$ws.Range("X2").Value = 0.01844135153902312
$ws.Range("Y2").Value = 0.08884395962632161
$ws.Range("Z2").Value = 0.6861817856681602
$ws.Range("AA2").Value = 0.1400038669760248
$ws.Range("AB2").Value = 0.01863158566370027
$ws.Range("AC2").Value = 0.1213722813123245
$ws.Range("AD2").Value = 84.06999999999999
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 84.06999999999999
$ws.Range("AG2").Value = 47.23999999999999
$ws.Range("AH2").Value = 0.1631666796056207
$ws.Range("AI2").Value = 0.5454486472458314
$ws.Range("AJ2").Value = 0.09874375535628434
$ws.Range("AK2").Value = 0.4027280477408354
$ws.Range("AL2").Value = 3.25
$ws.Range("AM2").Value = 2.561
$ws.Range("AN2").Value = 4.25240263024785
$ws.Range("AO2").Value = 5.833846153846153
$ws.Range("AP2").Value = 2.389479008598887
$ws.Range("AQ2").Value = 7.403358063256539

# --- Row 3: General Capital Limited (NZSE:GEN) - was NZX Limited ---
$ws.Range("B3").Value = "General Capital Limited (NZSE:GEN)"
$ws.Range("G3").Value = 0.5354609929078015
$ws.Range("H3").Value = 0.5354609929078015
$ws.Range("I3").Value = 0.7304964539007093
$ws.Range("J3").Value = 0.5025401136763744
$ws.Range("K3").Value = 0.097
$ws.Range("L3").Value = 0.03439716312056738
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 7.43
$ws.Range("V3").Value = 0.9537869062901155
$ws.Range("W3").Value = 0.0162751677852349
$ws.Range("X3").Value = 0.01705068102297429
$ws.Range("Y3").Value = -0.0007755132377393892
$ws.Range("Z3").Value = 1.2
$ws.Range("AA3").Value = 0.6030481364116492
$ws.Range("AB3").Value = 0.01705068102297429
$ws.Range("AC3").Value = 0.585997455388675
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -7.43
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -20.63888888888887
$ws.Range("AK3").Value = 6.19166666666667
$ws.Range("AL3").Value = 1.28
$ws.Range("AM3").Value = 1.28
$ws.Range("AN3").Value = 0
$ws.Range("AO3").Value = 1.609375
$ws.Range("AP3").Value = -3.589371980676328
$ws.Range("AQ3").Value = 1.609375
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("T3").ClearContents()

# --- Row 4: NZX Limited (NZSE:NZX) - was Blackwell Global Holdings Limited ---
$ws.Range("B4").Value = "NZX Limited (NZSE:NZX)"
$ws.Range("D4").Value = 0.0188
$ws.Range("E4").Value = -0.0645
$ws.Range("F4").Value = 0.107
$ws.Range("G4").Value = 0.3402061855670103
$ws.Range("H4").Value = 0.3402061855670103
$ws.Range("I4").Value = 0.3484536082474227
$ws.Range("J4").Value = 0.2488316151202749
$ws.Range("K4").Value = 11.2
$ws.Range("L4").Value = 0.2309278350515464
$ws.Range("M4").Value = 9.460000000000001
$ws.Range("N4").Value = 0.02412649834225963
$ws.Range("O4").Value = 0.8446428571428573
$ws.Range("P4").Value = 9.460000000000001
$ws.Range("Q4").Value = 0.02412649834225963
$ws.Range("R4").Value = 0.8446428571428573
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 13.8
$ws.Range("V4").Value = 0.03519510328997705
$ws.Range("W4").Value = 0.2685851318944844
$ws.Range("X4").Value = 0.01749045952692515
$ws.Range("Y4").Value = 0.2510946723675592
$ws.Range("Z4").Value = 1.125290023201856
$ws.Range("AA4").Value = 0.2800077339520495
$ws.Range("AB4").Value = 0.01722262595864622
$ws.Range("AC4").Value = 0.2627851079934033
$ws.Range("AD4").Value = 30.9
$ws.Range("AF4").Value = 30.9
$ws.Range("AG4").Value = 17.1
$ws.Range("AH4").Value = 0.07304964539007092
$ws.Range("AI4").Value = 0.4221311475409836
$ws.Range("AJ4").Value = 0.04178885630498533
$ws.Range("AK4").Value = 0.2878787878787879
$ws.Range("AL4").Value = 1.97
$ws.Range("AM4").Value = 1.281
$ws.Range("AN4").Value = 1.745762711864407
$ws.Range("AO4").Value = 8.578680203045685
$ws.Range("AP4").Value = 0.9661016949152541
$ws.Range("AQ4").Value = 13.1928181108509

# --- Row 5: Blackwell Global Holdings Limited (NZSE:BGI) - was Geneva Finance Limited ---
$ws.Range("B5").Value = "Blackwell Global Holdings Limited (NZSE:BGI)"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = -0.452
$ws.Range("L5").Value = -1.815261044176707
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("U5").Value = 1.5
$ws.Range("V5").Value = 0.3768844221105528
$ws.Range("W5").Value = -1.034324942791762
$ws.Range("X5").Value = 0.0193922435511211
$ws.Range("Y5").Value = -1.053717186342883
$ws.Range("Z5").Value = 0.190512624330528
$ws.Range("AA5").Value = 0
$ws.Range("AB5").Value = 0.02004054536875433
$ws.Range("AC5").Value = -0.02004054536875433
$ws.Range("AD5").Value = 1.67
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 1.67
$ws.Range("AG5").Value = 0.1699999999999999
$ws.Range("AH5").Value = 0.295575221238938
$ws.Range("AI5").Value = 0.8789473684210526
$ws.Range("AJ5").Value = 0.04096385542168673
$ws.Range("AK5").Value = 0.4249999999999999
$ws.Range("D5").ClearContents()
$ws.Range("T5").ClearContents()
$ws.Range("AN5").ClearContents()
$ws.Range("AP5").ClearContents()

# --- Row 6: Geneva Finance Limited (NZSE:GFL) - was General Capital Limited ---
$ws.Range("B6").Value = "Geneva Finance Limited (NZSE:GFL)"
$ws.Range("D6").Value = 0.25
$ws.Range("E6").Value = 0.091
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3.49
$ws.Range("L6").Value = 0.1994285714285714
$ws.Range("M6").Value = 1.46
$ws.Range("N6").Value = 0.05347985347985348
$ws.Range("O6").Value = 0.4183381088825214
$ws.Range("P6").Value = 1.46
$ws.Range("Q6").Value = 0.05347985347985348
$ws.Range("R6").Value = 0.4183381088825214
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 14.1
$ws.Range("V6").Value = 0.5164835164835164
$ws.Range("W6").Value = 0.1982954545454546
$ws.Range("X6").Value = 0.02757797964318992
$ws.Range("Y6").Value = 0.1707174749022646
$ws.Range("Z6").Value = 0.3246753246753246
$ws.Range("AA6").Value = 0
$ws.Range("AB6").Value = 0.02366164015557214
$ws.Range("AC6").Value = -0.02366164015557214
$ws.Range("AD6").Value = 51.5
$ws.Range("AF6").Value = 51.5
$ws.Range("AG6").Value = 37.4
$ws.Range("AH6").Value = 0.6535532994923858
$ws.Range("AI6").Value = 0.7074175824175825
$ws.Range("AJ6").Value = 0.5780525502318392
$ws.Range("AK6").Value = 0.6371379897785349
$ws.Range("AL6").Value = 0
$ws.Range("AM6").Value = 0
$ws.Range("AN6").ClearContents()
$ws.Range("AO6").ClearContents()
$ws.Range("AP6").ClearContents()
$ws.Range("AQ6").ClearContents()
